$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1342.6875
$ws.Range("I32").Value = 1296.6666
$ws.Range("K32").Value = 1296.6666
$ws.Range("M32").Value = -970.6666
$ws.Range("H33").Value = 257.2963
$ws.Range("I33").Value = 279.375
$ws.Range("J33").Value = 80.666664
$ws.Range("K33").Value = 279.375
$ws.Range("L33").Value = 80.666664
$ws.Range("M33").Value = -50.375
$ws.Range("N33").Value = -538.666664
$ws.Range("H74").Value = 3432
$ws.Range("I74").Value = 3630
$ws.Range("J74").Value = 3300
$ws.Range("K74").Value = 3630
$ws.Range("L74").Value = 3300
$ws.Range("M74").Value = -2694
$ws.Range("N74").Value = -5172
$ws.Range("H77").Value = 3432
$ws.Range("I77").Value = 3630
$ws.Range("J77").Value = 3300
$ws.Range("K77").Value = 18150
$ws.Range("L77").Value = 16500
$ws.Range("M77").Value = -13470
$ws.Range("N77").Value = -25860
$ws.Range("H132").Value = 2454.349
$ws.Range("I132").Value = 1929.2203
$ws.Range("J132").Value = 10200
$ws.Range("K132").Value = 5787.6609
$ws.Range("L132").Value = 30600
$ws.Range("M132").Value = -3257.6609
$ws.Range("N132").Value = -35660
$ws.Range("H137").Value = 5406270
$ws.Range("I137").Value = 809.9231
$ws.Range("J137").Value = 18182812
$ws.Range("K137").Value = 2429.7693
$ws.Range("L137").Value = 54548436
$ws.Range("M137").Value = 120.2307000000001
$ws.Range("N137").Value = -54553536
$ws.Range("H138").Value = 2239.362
$ws.Range("I138").Value = 1193.9722
$ws.Range("J138").Value = 3950
$ws.Range("K138").Value = 3581.9166
$ws.Range("L138").Value = 11850
$ws.Range("M138").Value = 1558.0834
$ws.Range("N138").Value = -22130

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 11629358
$ws.Range("I61").Value = 14707254
$ws.Range("J61").Value = 1753.2222
$ws.Range("K61").Value = 14707254
$ws.Range("L61").Value = 1753.2222
$ws.Range("M61").Value = -14707042
$ws.Range("N61").Value = -2177.2222
$ws.Range("H63").Value = 2750
$ws.Range("I63").Value = 2500
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 2500
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -1814
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 2750
$ws.Range("I66").Value = 2500
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 12500
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -9068
$ws.Range("N66").Value = -21864
$ws.Range("H132").Value = 7814921.5
$ws.Range("I132").Value = 11365792
$ws.Range("J132").Value = 3005.6
$ws.Range("K132").Value = 34097376
$ws.Range("L132").Value = 9016.799999999999
$ws.Range("M132").Value = -34094846
$ws.Range("N132").Value = -14076.8
$ws.Range("H136").Value = 11629358
$ws.Range("I136").Value = 14707254
$ws.Range("J136").Value = 1753.2222
$ws.Range("K136").Value = 44121762
$ws.Range("L136").Value = 5259.6666
$ws.Range("M136").Value = -44119212
$ws.Range("N136").Value = -10359.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 21741492
$ws.Range("I86").Value = 1908
$ws.Range("J86").Value = 38464250
$ws.Range("K86").Value = 1908
$ws.Range("L86").Value = 38464250
$ws.Range("M86").Value = -785
$ws.Range("N86").Value = -38466496
$ws.Range("H89").Value = 21741492
$ws.Range("I89").Value = 1908
$ws.Range("J89").Value = 38464250
$ws.Range("K89").Value = 9540
$ws.Range("L89").Value = 192321250
$ws.Range("M89").Value = -3924
$ws.Range("N89").Value = -192332482
$ws.Range("H134").Value = 5289.3335
$ws.Range("I134").Value = 5042.4
$ws.Range("J134").Value = 5412.8
$ws.Range("K134").Value = 15127.2
$ws.Range("L134").Value = 16238.4
$ws.Range("M134").Value = -12592.2
$ws.Range("N134").Value = -21308.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2728
$ws.Range("I58").Value = 987.1429000000001
$ws.Range("J58").Value = 3946.6
$ws.Range("K58").Value = 987.1429000000001
$ws.Range("L58").Value = 3946.6
$ws.Range("M58").Value = -784.1429000000001
$ws.Range("N58").Value = -4352.6
$ws.Range("H136").Value = 2728
$ws.Range("I136").Value = 987.1429000000001
$ws.Range("J136").Value = 3946.6
$ws.Range("K136").Value = 2961.4287
$ws.Range("L136").Value = 11839.8
$ws.Range("M136").Value = -411.4287000000004
$ws.Range("N136").Value = -16939.8
$ws.Range("H140").Value = 45608.43
$ws.Range("J140").Value = 45608.43
$ws.Range("L140").Value = 45608.43
$ws.Range("N140").Value = -55968.43

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 6836.4287
$ws.Range("I103").Value = 671.25
$ws.Range("J103").Value = 8287.058999999999
$ws.Range("K103").Value = 2013.75
$ws.Range("L103").Value = 24861.177
$ws.Range("M103").Value = -1134.75
$ws.Range("N103").Value = -26619.177
$ws.Range("H107").Value = 382.87097
$ws.Range("I107").Value = 185.66667
$ws.Range("J107").Value = 404
$ws.Range("K107").Value = 557.00001
$ws.Range("L107").Value = 1212
$ws.Range("M107").Value = 1362.99999
$ws.Range("N107").Value = -5052
$ws.Range("H118").Value = 1203.7273
$ws.Range("I118").Value = 900
$ws.Range("J118").Value = 1234.1
$ws.Range("K118").Value = 2700
$ws.Range("L118").Value = 3702.3
$ws.Range("M118").Value = -1457
$ws.Range("N118").Value = -6188.299999999999
$ws.Range("H134").Value = 3297.2812
$ws.Range("I134").Value = 1868.7727
$ws.Range("J134").Value = 6440
$ws.Range("K134").Value = 5606.3181
$ws.Range("L134").Value = 19320
$ws.Range("M134").Value = -536.3181000000004
$ws.Range("N134").Value = -29460
$ws.Range("H139").Value = 2383.6667
$ws.Range("I139").Value = 1462.1538
$ws.Range("J139").Value = 4779.6
$ws.Range("K139").Value = 4386.4614
$ws.Range("L139").Value = 14338.8
$ws.Range("M139").Value = 753.5385999999999
$ws.Range("N139").Value = -24618.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 33338144
$ws.Range("I80").Value = 66668788
$ws.Range("J80").Value = 7500
$ws.Range("K80").Value = 66668788
$ws.Range("L80").Value = 7500
$ws.Range("M80").Value = -66667790
$ws.Range("N80").Value = -9496
$ws.Range("H83").Value = 33338144
$ws.Range("I83").Value = 66668788
$ws.Range("J83").Value = 7500
$ws.Range("K83").Value = 333343940
$ws.Range("L83").Value = 37500
$ws.Range("M83").Value = -333338948
$ws.Range("N83").Value = -47484
$ws.Range("H113").Value = 168163.67
$ws.Range("I113").Value = 168163.67
$ws.Range("K113").Value = 168163.67
$ws.Range("M113").Value = -165993.67
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 3480.8298
$ws.Range("I132").Value = 2602.6365
$ws.Range("J132").Value = 5550.857
$ws.Range("K132").Value = 7807.9095
$ws.Range("L132").Value = 16652.571
$ws.Range("M132").Value = -5277.9095
$ws.Range("N132").Value = -21712.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1093.7727
$ws.Range("I22").Value = 508.75
$ws.Range("J22").Value = 1428.0714
$ws.Range("K22").Value = 508.75
$ws.Range("L22").Value = 1428.0714
$ws.Range("M22").Value = -213.75
$ws.Range("N22").Value = -2018.0714
$ws.Range("H27").Value = 1093.7727
$ws.Range("I27").Value = 508.75
$ws.Range("J27").Value = 1428.0714
$ws.Range("K27").Value = 508.75
$ws.Range("L27").Value = 1428.0714
$ws.Range("M27").Value = -401.75
$ws.Range("N27").Value = -1642.0714
$ws.Range("H55").Value = 989.1111
$ws.Range("I55").Value = 414.57144
$ws.Range("J55").Value = 3000
$ws.Range("K55").Value = 414.57144
$ws.Range("L55").Value = 3000
$ws.Range("M55").Value = -241.57144
$ws.Range("N55").Value = -3346
$ws.Range("H82").Value = 2440.2
$ws.Range("I82").Value = 2167
$ws.Range("J82").Value = 2850
$ws.Range("K82").Value = 2167
$ws.Range("L82").Value = 2850
$ws.Range("M82").Value = -1806
$ws.Range("N82").Value = -3572
$ws.Range("H85").Value = 2440.2
$ws.Range("I85").Value = 2167
$ws.Range("J85").Value = 2850
$ws.Range("K85").Value = 2167
$ws.Range("L85").Value = 2850
$ws.Range("M85").Value = -919
$ws.Range("N85").Value = -5346
$ws.Range("H140").Value = 52582.25
$ws.Range("J140").Value = 52582.25
$ws.Range("L140").Value = 52582.25
$ws.Range("N140").Value = -62942.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1177.8422
$ws.Range("I100").Value = 1013
$ws.Range("J100").Value = 1297.7273
$ws.Range("K100").Value = 2026
$ws.Range("L100").Value = 2595.4546
$ws.Range("M100").Value = -1485
$ws.Range("N100").Value = -3677.4546
